$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new news item
$ws.Range("A2").Value = "Prabowo Rapat Bareng Menteri di Kertanegara Minggu Malam, Bahas MBG hingga Migas"
$ws.Range("B2").Value = "2025-09-29T01:01:30+07:00"
$ws.Range("D2").Value = "https://www.liputan6.com/news/read/6170969/prabowo-rapat-bareng-menteri-di-kertanegara-minggu-malam-bahas-mbg-hingga-migas"

# Delete rows 3 through 6 entirely, leaving only header + row 2
$ws.Range("A3:E6").EntireRow.Delete()

$wb.Save()
